$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new log entry as row 47 (4th July, 8:30pm - 9:30pm, Code) ---
# Copy the formatting of the previous entry row (46) onto the new row (47)
# first, so number formats / alignment / borders / wrap-text all match the
# rest of the log table.
$ws.Range("A46:G46").Copy()
$ws.Range("A47:G47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 44746
$ws.Range("C47").Value = 0.85416666666666663
$ws.Range("D47").Value = 0.89583333333333337
$ws.Range("E47").Formula = "=D47-C47"
$ws.Range("F47").Value = "Code"
$ws.Range("G47").Value = "1. Label output formats `n2. Loss function presentation`n3. Added weight files, videos to nb, drive"
$ws.Range("G47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 45

# --- Move the view / selection to reflect the newly added row ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 41
$aw.ScrollColumn = 1
$null = $ws.Range("G48").Select()
